# Update "paises.xlsx" (COVID country stats) and "provincias Spain" style
# refresh: new scrape snapshot + timestamp, which re-sorts a few rows whose
# totals become very close to their neighbours (hence a handful of country
# labels trade places while keeping the same descending order by "Casos
# totales").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Septiembre de 2020 a las 13:56"

# --- Estados Unidos (row 4) --------------------------------------------
$ws.Cells.Item(4, 2).Value = 7005686
$ws.Cells.Item(4, 3).Value = 918
$ws.Cells.Item(4, 4).Value = 4250497
$ws.Cells.Item(4, 5).Value = 2551067
$ws.Cells.Item(4, 7).Value = 4
$ws.Cells.Item(4, 8).Value = 204122

# --- India (row 5) -------------------------------------------------------
$ws.Cells.Item(5, 2).Value = 5491410
$ws.Cells.Item(5, 3).Value = 5798
$ws.Cells.Item(5, 5).Value = 1007078
$ws.Cells.Item(5, 7).Value = 24
$ws.Cells.Item(5, 8).Value = 87933

# --- Iran (row 16) ---------------------------------------------------
$ws.Cells.Item(16, 2).Value = 425481
$ws.Cells.Item(16, 3).Value = 3341
$ws.Cells.Item(16, 4).Value = 361523
$ws.Cells.Item(16, 5).Value = 39480
$ws.Cells.Item(16, 7).Value = 177
$ws.Cells.Item(16, 8).Value = 24478

# --- Alemania (row 25) ------------------------------------------------
$ws.Cells.Item(25, 2).Value = 273793
$ws.Cells.Item(25, 3).Value = 316
$ws.Cells.Item(25, 5).Value = 20322
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 9471

# --- Catar (row 32) -----------------------------------------------------
$ws.Cells.Item(32, 2).Value = 123604
$ws.Cells.Item(32, 3).Value = 228
$ws.Cells.Item(32, 4).Value = 120540
$ws.Cells.Item(32, 5).Value = 2853
$ws.Cells.Item(32, 7).Value = 1
$ws.Cells.Item(32, 8).Value = 211

# --- Paises Bajos / Oman trade places (rows 41-42) ----------------------
$ws.Cells.Item(41, 1).Value = "Oman"
$ws.Cells.Item(41, 2).Value = 94051
$ws.Cells.Item(41, 3).Value = 576
$ws.Cells.Item(41, 4).Value = 85781
$ws.Cells.Item(41, 5).Value = 7417
$ws.Cells.Item(41, 7).Value = 7
$ws.Cells.Item(41, 8).Value = 853

$ws.Cells.Item(42, 1).Value = "Paises Bajos"
$ws.Cells.Item(42, 2).Value = 93778
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 8).Value = 6279

# --- Barein / Nepal trade places (rows 54-55) ---------------------------
$ws.Cells.Item(54, 1).Value = "Nepal"
$ws.Cells.Item(54, 2).Value = 65276
$ws.Cells.Item(54, 3).Value = 1154
$ws.Cells.Item(54, 4).Value = 47238
$ws.Cells.Item(54, 5).Value = 17611
$ws.Cells.Item(54, 7).Value = 16
$ws.Cells.Item(54, 8).Value = 427

$ws.Cells.Item(55, 1).Value = "Barein"
$ws.Cells.Item(55, 2).Value = 65039
$ws.Cells.Item(55, 4).Value = 57950
$ws.Cells.Item(55, 5).Value = 6866
$ws.Cells.Item(55, 7).Value = 2
$ws.Cells.Item(55, 8).Value = 223

# --- Argelia / Chequia / Suiza reshuffle (rows 60-62) -------------------
$ws.Cells.Item(60, 1).Value = "Suiza"
$ws.Cells.Item(60, 2).Value = 50378
$ws.Cells.Item(60, 3).Value = 1095
$ws.Cells.Item(60, 4).Value = 40500
$ws.Cells.Item(60, 5).Value = 7830
$ws.Cells.Item(60, 7).Value = 3
$ws.Cells.Item(60, 8).Value = 2048

$ws.Cells.Item(61, 1).Value = "Argelia"
$ws.Cells.Item(61, 2).Value = 49826
$ws.Cells.Item(61, 4).Value = 35047
$ws.Cells.Item(61, 5).Value = 13107
$ws.Cells.Item(61, 8).Value = 1672

$ws.Cells.Item(62, 1).Value = "Chequia"
$ws.Cells.Item(62, 2).Value = 49290
$ws.Cells.Item(62, 4).Value = 24755
$ws.Cells.Item(62, 5).Value = 24032
$ws.Cells.Item(62, 8).Value = 503

# --- Madagascar (row 87) -------------------------------------------------
$ws.Cells.Item(87, 2).Value = 16073
$ws.Cells.Item(87, 3).Value = 20
$ws.Cells.Item(87, 4).Value = 14682
$ws.Cells.Item(87, 5).Value = 1166
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 225

# --- Senegal (row 90) -----------------------------------------------------
$ws.Cells.Item(90, 2).Value = 14738
$ws.Cells.Item(90, 3).Value = 24
$ws.Cells.Item(90, 4).Value = 11458
$ws.Cells.Item(90, 5).Value = 2978

# --- Timor Oriental / Santa Lucia trade places (rows 204-205, tied totals)
$ws.Cells.Item(204, 1).Value = "Santa Lucia"
$ws.Cells.Item(205, 1).Value = "Timor Oriental"

# --- Islas Malvinas / Montserrat trade places (rows 214-215) -------------
$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 8).Value = 1

$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0
